$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new review row reuses row 2's layout/formatting (same appid+keyword,
# same column styles) so clone row 2 down into row 3 first.
$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))

# Fill in the new review's own values.
$ws.Range("C3").Value = "jorjkluni03@gmail.com"
$ws.Range("D3").Value = "vikicrestina@gmail.com"
$ws.Range("E3").Value = "27/5/2019 15:58"
$ws.Range("F3").Value = "One of my favorites car games in google play. Recommended!"

# Email columns are mailto hyperlinks, same as row 2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")

# Adding a hyperlink stamps Excel's built-in blue/underline "Hyperlink" style
# on the cell; restore the original (copied-from-row-2) cell formatting so
# C3/D3 keep looking like the rest of the table.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("F3").Select() | Out-Null
